$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 13.103
$ws.Range("E8").Value = 13.419
$ws.Range("D12").Value = -8.130000000000001
$ws.Range("E12").Value = 13.004
$ws.Range("E14").Value = 13.072
$ws.Range("E22").Value = 13.107
